$d = $word.ActiveDocument

$replacements = @(
    @("muundo. At their disposal there are", "muundo. Ovyo wao wapo"),
    @("some snails a big metal ring and a long", "konokono wengine pete kubwa ya chuma na ndefu"),
    @("string.", "kamba."),
    @("Explain how the team can manage to use", "Eleza jinsi timu inaweza kusimamia matumizi"),
    @("the materials to tell approximately the", "nyenzo za kusema takriban"),
    @("ideal location of the airport. Imagine", "eneo bora la uwanja wa ndege. Fikiria"),
    @("that the cities are placed at the", "kwamba miji imewekwa kwenye"),
    @("vertices of a triangle which is", "vipeo vya pembetatu ambayo ni"),
    @("obviously reproduced in scale as", "kwa hakika imetolewa tena kwa kiwango kama"),
    @("shown in figure. This is one possible", "inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana"),
    @("setting the rope starts from one nail,", "kuweka kamba huanza kutoka msumari mmoja,"),
    @("goes inside the ring, goes around the", "huenda ndani ya pete, huzunguka"),
    @("other nail, the third nail, inside the", "msumari mwingine, msumari wa tatu, ndani ya"),
    @("ring again and now you can just pull the", "pete tena na sasa unaweza kuvuta tu"),
    @("rope in order to find the point that", "kamba ili kupata uhakika huo"),
    @("you're looking for. In order to reach the", "unatafuta. Ili kufikia"),
    @("point, we have to move the rope a bit", "uhakika, tunapaswa kusonga kamba kidogo"),
    @("because there is some ", "kwa sababu kuna "),
    @("resistance", "upinzani"),
    @(" caused", " uliosababishwa"),
    @("by the materials that we are using but", "kwa nyenzo ambazo tunatumia lakini"),
    @("after a while you'll reach a position from", "baada ya muda utafikia nafasi kutoka"),
    @("which the ring doesn't move anymore,", "ambayo pete haisogei tena,"),
    @("which is more or less this one. And as", "ambayo ni zaidi au chini ya hii. Na kama"),
    @("between the ring and the nails are", "kati ya pete na misumari ni"),
    @("placed more or less 120 degrees from one", "kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja"),
    @("another which is 1/3 of a circumference,", "nyingine ambayo ni 1/3 ya mduara,"),
    @("and that's the point that we're looking", "na hiyo ndiyo hatua tunayoiangalia"),
    @("for: the minimum distance between the", "kwa: umbali wa chini kati ya"),
    @("nails and the airport when you sum it", "misumari na uwanja wa ndege unapojumlisha"),
    @("ogether", "pamoja"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
